# Adds a new data row (row 3) to the "mooring", "fender" and "lngc"
# worksheets of the pretension analysis summary workbook, capturing the
# results for the new starboard (SB) vessel-statics run:
#   fsts_l015_hwl_125km3_l100_sb_vessel_statics_6dof
#
# Row 2 in each sheet already holds the equivalent "pb" (port) run; row 3
# appends the matching "sb" (starboard) run directly below it.

$wb = $excel.ActiveWorkbook

$simFile = "D:/github/digitalmodel/specs/modules/orcaflex/mooring-tension-iteration/go-by/.sim/fsts_l015_hwl_125km3_l100_sb_vessel_statics_6dof.sim"
$stem    = "fsts_l015_hwl_125km3_l100_sb_vessel_statics_6dof"

# ---------------------------------------------------------------------
# Sheet "mooring" -> row 3 (line tensions, columns A-E then H-W)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("mooring")

$ws.Range("A3").Value = $simFile
$ws.Range("B3").Value = $stem
$ws.Range("C3").Value = "InStaticState"
$ws.Range("D3").Value = -10
$ws.Range("E3").Value = "inf"

$ws.Range("H3").Value = 99.441947
$ws.Range("I3").Value = 99.43378
$ws.Range("J3").Value = 99.401663
$ws.Range("K3").Value = 76.19195000000001
$ws.Range("L3").Value = 92.410622
$ws.Range("M3").Value = 60.658374
$ws.Range("N3").Value = 120.644396
$ws.Range("O3").Value = 120.650893
$ws.Range("P3").Value = 119.326597
$ws.Range("Q3").Value = 119.319975
$ws.Range("R3").Value = 119.342562
$ws.Range("S3").Value = 120.544677
$ws.Range("T3").Value = 120.297567
$ws.Range("U3").Value = 120.590698
$ws.Range("V3").Value = 120.433381
$ws.Range("W3").Value = 120.419642

# ---------------------------------------------------------------------
# Sheet "fender" -> row 3 (fender contact forces, columns A-E, H-M, O)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("fender")

$ws.Range("A3").Value = $simFile
$ws.Range("B3").Value = $stem
$ws.Range("C3").Value = "InStaticState"
$ws.Range("D3").Value = -10
$ws.Range("E3").Value = "inf"

$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 249.473066
$ws.Range("J3").Value = 228.6993
$ws.Range("K3").Value = 142.165709
$ws.Range("L3").Value = 117.937571
$ws.Range("M3").Value = 23.590793
$ws.Range("O3").Value = 0

# ---------------------------------------------------------------------
# Sheet "lngc" -> row 3 (LNGC offsets/rotations, columns A-E, H-M)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("lngc")

$ws.Range("A3").Value = $simFile
$ws.Range("B3").Value = $stem
$ws.Range("C3").Value = "InStaticState"
$ws.Range("D3").Value = -10
$ws.Range("E3").Value = "inf"

$ws.Range("H3").Value = 76.88400300000001
$ws.Range("I3").Value = -51.688469
$ws.Range("J3").Value = -11.409393
$ws.Range("K3").Value = -0.031814
$ws.Range("L3").Value = -0.000557
$ws.Range("M3").Value = -179.98671
